$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-16 18:18:38"
$ws.Range("I2").Value = "19.1 mm"
$ws.Range("E3").Value = "2026-02-16 18:18:40"
$ws.Range("G3").Value = "234 cm"
$ws.Range("I3").Value = "9.8 mm"
$ws.Range("L3").Value = "68.8 km/h - 246º 17:50 TU"
$ws.Range("O3").Value = "-0.8 °C"
$ws.Range("E4").Value = "2026-02-16 18:18:42"
$ws.Range("J4").Value = "1012.4 hPa"
$ws.Range("E5").Value = "2026-02-16 18:18:45"
$ws.Range("I5").Value = "22.8 mm"
$ws.Range("L5").Value = "43.9 km/h - 327º 17:30 TU"
$ws.Range("N5").Value = "-1.5 °C 17:41 TU"
$ws.Range("O5").Value = "-0.6 °C"
$ws.Range("E6").Value = "2026-02-16 18:18:47"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "68%"
$ws.Range("J6").Value = "1012.5 hPa"
$ws.Range("O6").Value = "11.7 °C"
$ws.Range("E7").Value = "2026-02-16 18:18:50"
$ws.Range("J7").Value = "1013.5 hPa"
$ws.Range("O7").Value = "16.3 °C"
$ws.Range("E8").Value = "2026-02-16 18:18:52"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "65%"
$ws.Range("J8").Value = "1013.1 hPa"
$ws.Range("O8").Value = "12.5 °C"
$ws.Range("E9").Value = "2026-02-16 18:18:55"
$ws.Range("O9").Value = "11.1 °C"
$ws.Range("E10").Value = "2026-02-16 18:18:57"
$ws.Range("O10").Value = "10.9 °C"
$ws.Range("E11").Value = "2026-02-16 18:19:00"
$ws.Range("O11").Value = "6.6 °C"
$ws.Range("E12").Value = "2026-02-16 18:19:02"
$ws.Range("O12").Value = "10.5 °C"
$ws.Range("E13").Value = "2026-02-16 18:19:05"
$ws.Range("J13").Value = "1014.9 hPa"
$ws.Range("O13").Value = "5.6 °C"
$ws.Range("E14").Value = "2026-02-16 18:19:07"
$ws.Range("O14").Value = "16.0 °C"
$ws.Range("E15").Value = "2026-02-16 18:19:10"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "69%"
$ws.Range("O15").Value = "11.2 °C"
$ws.Range("E16").Value = "2026-02-16 18:19:12"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "80%"
$ws.Range("O16").Value = "0.0 °C"
$ws.Range("E17").Value = "2026-02-16 18:19:15"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "69%"
$ws.Range("N17").Value = "4.1 °C 17:53 TU"
$ws.Range("E18").Value = "2026-02-16 18:19:17"
$ws.Range("J18").Value = "1012.8 hPa"
$ws.Range("O18").Value = "10.9 °C"
$ws.Range("E19").Value = "2026-02-16 18:19:20"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "82%"
$ws.Range("O19").Value = "6.9 °C"
$ws.Range("E20").Value = "2026-02-16 18:19:22"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "98%"
$ws.Range("I20").Value = "0.2 mm"
$ws.Range("O20").Value = "-0.6 °C"
$ws.Range("E21").Value = "2026-02-16 18:19:25"
$ws.Range("J21").Value = "1014.4 hPa"
$ws.Range("O21").Value = "8.3 °C"
$ws.Range("E22").Value = "2026-02-16 18:19:27"
$ws.Range("E23").Value = "2026-02-16 18:19:30"
$ws.Range("I23").Value = "12.9 mm"
$ws.Range("L23").Value = "67.7 km/h - 261º 17:42 TU"
$ws.Range("O23").Value = "-0.6 °C"
$ws.Range("E24").Value = "2026-02-16 18:19:32"
$ws.Range("J24").Value = "1016.7 hPa"
$ws.Range("O24").Value = "12.9 °C"
$ws.Range("E25").Value = "2026-02-16 18:19:35"
$ws.Range("I25").Value = "5.3 mm"
$ws.Range("E26").Value = "2026-02-16 18:19:37"
$ws.Range("E27").Value = "2026-02-16 18:19:40"
$ws.Range("L27").Value = "47.2 km/h - 252º 17:54 TU"
$ws.Range("E28").Value = "2026-02-16 18:19:42"
$ws.Range("J28").Value = "1012.9 hPa"
$ws.Range("O28").Value = "9.5 °C"
$ws.Range("E29").Value = "2026-02-16 18:19:45"
$ws.Range("O29").Value = "10.8 °C"
$ws.Range("E30").Value = "2026-02-16 18:19:47"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "68%"
$ws.Range("O30").Value = "11.9 °C"
$ws.Range("E31").Value = "2026-02-16 18:19:50"
$ws.Range("E32").Value = "2026-02-16 18:19:52"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "80%"
$ws.Range("O32").Value = "8.5 °C"
$ws.Range("E33").Value = "2026-02-16 18:19:55"
$ws.Range("E34").Value = "2026-02-16 18:19:57"
$ws.Range("O34").Value = "3.6 °C"
$ws.Range("E35").Value = "2026-02-16 18:20:00"
$ws.Range("J35").Value = "1016.9 hPa"
$ws.Range("O35").Value = "9.5 °C"
$ws.Range("E36").Value = "2026-02-16 18:20:02"
$ws.Range("O36").Value = "11.6 °C"
$ws.Range("E37").Value = "2026-02-16 18:20:05"
$ws.Range("O37").Value = "6.5 °C"
$ws.Range("E38").Value = "2026-02-16 18:20:07"
$ws.Range("O38").Value = "11.9 °C"
$ws.Range("E39").Value = "2026-02-16 18:20:10"
$ws.Range("G39").Value = "56 cm"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "77%"
$ws.Range("I39").Value = "2.4 mm"
$ws.Range("K39").Value = "11.1 MJ/m2"
$ws.Range("L39").Value = "71.3 km/h - 290º 16:11 TU"
$ws.Range("N39").Value = "-0.7 °C 17:53 TU"
$ws.Range("E40").Value = "2026-02-16 18:20:12"
$ws.Range("J40").Value = "1016.6 hPa"
$ws.Range("O40").Value = "6.8 °C"
$ws.Range("E41").Value = "2026-02-16 18:20:15"
$ws.Range("E42").Value = "2026-02-16 18:20:17"
$ws.Range("O42").Value = "11.2 °C"
$ws.Range("E43").Value = "2026-02-16 18:20:19"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "80%"
$ws.Range("O43").Value = "8.3 °C"
$ws.Range("E44").Value = "2026-02-16 18:20:22"
$ws.Range("N44").Value = "-1.9 °C 17:37 TU"
$ws.Range("O44").Value = "0.0 °C"
$ws.Range("E45").Value = "2026-02-16 18:20:24"
$ws.Range("I45").Value = "15.8 mm"
$ws.Range("E46").Value = "2026-02-16 18:20:27"
$ws.Range("O46").Value = "16.0 °C"
